$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.489.42"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "3.120.12"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'238.35"
$ws.Range("E5").Value = "  -2.82%  "
$ws.Range("D6").Value = "'615.48"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "'1.11"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("D8").Value = "'0.393"
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +13.29%  "
$ws.Range("D11").Value = "3.116.07"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "'0.199"
$ws.Range("E12").Value = "  -2.47%  "
$ws.Range("D13").Value = "'0.0000246"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("D14").Value = "'35.33"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").Value = "93.173.34"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("E16").Value = "  -3.14%  "
$ws.Range("D17").Value = "3.699.68"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "3.121.16"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("D20").Value = "'14.87"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("E21").Value = "  +3.98%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'443.79"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("B23").Value = "PEPE"
$ws.Range("C23").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D23").Value = "'0.0000201"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").Value = "'9.12"
$ws.Range("D25").Value = "'8.24"
$ws.Range("E25").Value = "  +4.93%  "
$ws.Range("D26").Value = "'5.79"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").Value = "'12.89"
$ws.Range("E27").Value = "  +10.10%  "
$ws.Range("D28").Value = "'85.84"
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("D30").Value = "'0.183"
$ws.Range("E30").Value = "  +9.71%  "
$ws.Range("D31").Value = "'0.239"
$ws.Range("E31").Value = "  +3.46%  "
$ws.Range("E32").Value = "  -11.31%  "
$ws.Range("D33").Value = "'9.29"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("E34").Value = "  +2.58%  "
$ws.Range("D35").Value = "'8.05"
$ws.Range("E35").Value = "  +4.02%  "
$ws.Range("E36").Value = "  -8.53%  "
$ws.Range("D37").Value = "'26.02"
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("B38").Value = "MantraDAO"
$ws.Range("C38").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D38").Value = "'3.96"
$ws.Range("E38").Value = "  -4.84%  "
$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D39").Value = "'1.91"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").Value = "'0.451"
$ws.Range("E40").Value = "  +2.57%  "
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").Value = "'479.13"
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("D43").Value = "'24.01"
$ws.Range("E43").Value = "  +8.14%  "
$ws.Range("D44").Value = "'3.33"
$ws.Range("E44").Value = "  -3.09%  "
$ws.Range("D46").Value = "'159.12"
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D47").Value = "'0.705"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").Value = "'4.41"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "'44.01"
$ws.Range("E51").Value = "  -0.33%  "
